$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27784882
$ws.Range("I32").Value = 7583.5864
$ws.Range("J32").Value = 142862260
$ws.Range("K32").Value = 7583.5864
$ws.Range("L32").Value = 142862260
$ws.Range("M32").Value = -7296.5864
$ws.Range("N32").Value = -142862834

$ws.Range("H63").Value = 2491.973
$ws.Range("J63").Value = 2626.923
$ws.Range("L63").Value = 2626.923
$ws.Range("N63").Value = -3998.923

$ws.Range("H66").Value = 2491.973
$ws.Range("J66").Value = 2626.923
$ws.Range("L66").Value = 13134.615
$ws.Range("N66").Value = -19998.615

$ws.Range("H80").Value = 20950
$ws.Range("J80").Value = 20950
$ws.Range("L80").Value = 20950
$ws.Range("N80").Value = -22946

$ws.Range("H83").Value = 20950
$ws.Range("J83").Value = 20950
$ws.Range("L83").Value = 62850
$ws.Range("N83").Value = -72834

$ws.Range("H122").Value = 1865.75
$ws.Range("I122").Value = 1728
$ws.Range("J122").Value = 2003.5
$ws.Range("K122").Value = 5184
$ws.Range("L122").Value = 6010.5
$ws.Range("M122").Value = -2734
$ws.Range("N122").Value = -10910.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2323.08
$ws.Range("I20").Value = 2672.4666
$ws.Range("J20").Value = 1799
$ws.Range("K20").Value = 2672.4666
$ws.Range("L20").Value = 1799
$ws.Range("M20").Value = -2425.4666
$ws.Range("N20").Value = -2293

$ws.Range("H82").Value = 10977.25
$ws.Range("I82").Value = 9211.4
$ws.Range("J82").Value = 19806.5
$ws.Range("K82").Value = 9211.4
$ws.Range("L82").Value = 19806.5
$ws.Range("M82").Value = -8828.4
$ws.Range("N82").Value = -20572.5

$ws.Range("H85").Value = 10977.25
$ws.Range("I85").Value = 9211.4
$ws.Range("J85").Value = 19806.5
$ws.Range("K85").Value = 9211.4
$ws.Range("L85").Value = 19806.5
$ws.Range("M85").Value = -7885.4
$ws.Range("N85").Value = -22458.5

$ws.Range("H86").Value = 2909689
$ws.Range("I86").Value = 2966.6667
$ws.Range("K86").Value = 2966.6667
$ws.Range("M86").Value = -1843.6667

$ws.Range("H89").Value = 2909689
$ws.Range("I89").Value = 2966.6667
$ws.Range("K89").Value = 14833.3335
$ws.Range("M89").Value = -9217.333500000001

$ws.Range("H134").Value = 2472312
$ws.Range("I134").Value = 831.6389
$ws.Range("J134").Value = 12358234
$ws.Range("K134").Value = 2494.9167
$ws.Range("L134").Value = 37074702
$ws.Range("M134").Value = 40.08329999999978
$ws.Range("N134").Value = -37079772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1933.3334
$ws.Range("J15").Value = 1933.3334
$ws.Range("L15").Value = 1933.3334
$ws.Range("N15").Value = -2273.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 206
$ws.Range("I14").Value = 213.33333
$ws.Range("J14").Value = 198.66667
$ws.Range("K14").Value = 213.33333
$ws.Range("L14").Value = 198.66667
$ws.Range("M14").Value = -45.33332999999999
$ws.Range("N14").Value = -534.6666700000001

$ws.Range("H21").Value = 25980
$ws.Range("J21").Value = 25980
$ws.Range("L21").Value = 25980
$ws.Range("N21").Value = -26326

$ws.Range("H30").Value = 25980
$ws.Range("J30").Value = 25980
$ws.Range("L30").Value = 25980
$ws.Range("N30").Value = -26190

$ws.Range("H57").Value = 15933.333
$ws.Range("J57").Value = 15933.333
$ws.Range("L57").Value = 15933.333
$ws.Range("N57").Value = -17573.333

$ws.Range("H62").Value = 8888
$ws.Range("J62").Value = 8888
$ws.Range("L62").Value = 8888
$ws.Range("N62").Value = -10260

$ws.Range("H65").Value = 8888
$ws.Range("J65").Value = 8888
$ws.Range("L65").Value = 26664
$ws.Range("N65").Value = -33528

$ws.Range("H69").Value = 11665.667
$ws.Range("J69").Value = 11665.667
$ws.Range("L69").Value = 11665.667
$ws.Range("N69").Value = -13163.667

$ws.Range("H70").Value = 6454.298
$ws.Range("I70").Value = 6828.825
$ws.Range("J70").Value = 4314.143
$ws.Range("K70").Value = 6828.825
$ws.Range("L70").Value = 4314.143
$ws.Range("M70").Value = -6558.825
$ws.Range("N70").Value = -4854.143

$ws.Range("H72").Value = 11665.667
$ws.Range("J72").Value = 11665.667
$ws.Range("L72").Value = 34997.001
$ws.Range("N72").Value = -42485.001

$ws.Range("H73").Value = 6454.298
$ws.Range("I73").Value = 6828.825
$ws.Range("J73").Value = 4314.143
$ws.Range("K73").Value = 6828.825
$ws.Range("L73").Value = 4314.143
$ws.Range("M73").Value = -5892.825
$ws.Range("N73").Value = -6186.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 700
$ws.Range("J3").Value = 700
$ws.Range("L3").Value = 700
$ws.Range("N3").Value = -924

$ws.Range("H14").Value = 9900
$ws.Range("J14").Value = 9900
$ws.Range("L14").Value = 9900
$ws.Range("N14").Value = -10244

$ws.Range("H15").Value = 700
$ws.Range("J15").Value = 700
$ws.Range("L15").Value = 700
$ws.Range("N15").Value = -1040

$ws.Range("H132").Value = 28578388
$ws.Range("I132").Value = 45715656
$ws.Range("J132").Value = 16273.4
$ws.Range("K132").Value = 137146968
$ws.Range("L132").Value = 48820.2
$ws.Range("M132").Value = -137144438
$ws.Range("N132").Value = -53880.2

$ws.Range("H136").Value = 164838000
$ws.Range("I136").Value = 114289110
$ws.Range("J136").Value = 333334270
$ws.Range("K136").Value = 342867330
$ws.Range("L136").Value = 1000002810
$ws.Range("M136").Value = -342864780
$ws.Range("N136").Value = -1000007910

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 20666.5
$ws.Range("I14").Value = 8000
$ws.Range("J14").Value = 33333
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 33333
$ws.Range("M14").Value = -7832
$ws.Range("N14").Value = -33669

$ws.Range("H20").Value = 40005.5
$ws.Range("J20").Value = 40005.5
$ws.Range("L20").Value = 40005.5
$ws.Range("N20").Value = -40485.5

$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 4000
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 4000
$ws.Range("M26").Value = -707
$ws.Range("N26").Value = -4586

$ws.Range("H32").Value = 1495
$ws.Range("I32").Value = 1495
$ws.Range("K32").Value = 1495
$ws.Range("M32").Value = -1178

$ws.Range("H52").Value = 3540
$ws.Range("I52").Value = 1600
$ws.Range("J52").Value = 4833.3335
$ws.Range("K52").Value = 1600
$ws.Range("L52").Value = 4833.3335
$ws.Range("M52").Value = -1374
$ws.Range("N52").Value = -5285.3335

$ws.Range("H98").Value = 39800
$ws.Range("J98").Value = 39800
$ws.Range("L98").Value = 39800
$ws.Range("N98").Value = -45790

$ws.Range("H132").Value = 39783.832
$ws.Range("I132").Value = 88201
$ws.Range("J132").Value = 7505.722
$ws.Range("K132").Value = 264603
$ws.Range("L132").Value = 22517.166
$ws.Range("M132").Value = -262073
$ws.Range("N132").Value = -27577.166

$ws.Range("H136").Value = 2200.325
$ws.Range("I136").Value = 1569.2
$ws.Range("J136").Value = 2579
$ws.Range("K136").Value = 4707.6
$ws.Range("L136").Value = 7737
$ws.Range("M136").Value = -2157.6
$ws.Range("N136").Value = -12837
